$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: nrSessao 1 -> 3 ; photovoice "" -> "NAO" ; apresentacao_photovoice "" -> "NAO"
$ws.Range("N3").Value = 3
$ws.Range("O3").Value = "NAO"
$ws.Range("P3").Value = "NAO"

# Row 4: nrSessao 1 -> 5 ; photovoice "" -> "NAO" ; apresentacao_photovoice "" -> "NAO"
$ws.Range("N4").Value = 5
$ws.Range("O4").Value = "NAO"
$ws.Range("P4").Value = "NAO"

# Row 5: nrSessao 5 -> 6 (photovoice / apresentacao_photovoice unchanged, already "NAO")
$ws.Range("N5").Value = 6
